# Compare_sightlines.xlsx edit
# "Modified interpolation to reflect resolution"
#
# Adds a "No. Observations co added" column to the Table2 table on the
# "Parallel_185859" sheet (sheet1), populates its values, nudges the big
# comparison-plot picture one column to the right to make room for it,
# sets the sheet to print in portrait orientation, and updates the
# active-sheet/selection state (Parallel_185859 / L13 becomes active,
# instead of Parallel_all).

$wb = $excel.ActiveWorkbook

$wsParallel      = $wb.Worksheets.Item("Parallel_185859")
$wsPerpendicular = $wb.Worksheets.Item("Perpendicular_185859")

# --- 1. Extend Table2 with the new "No. Observations co added" column ---
$table = $wsParallel.ListObjects.Item("Table2")
$newColumn = $table.ListColumns.Add()

# Setting the header cell's value (rather than the ListColumn.Name
# property) is what actually persists the new column name.
$table.HeaderRowRange.Item(1, 4).Value = "No. Observations co added"

# --- 2. Fill in the per-sightline observation counts ---
$counts = @{
    6  = 2
    7  = 1
    8  = 1
    9  = 3
    10 = 2
    11 = 1
    12 = 5
    13 = 3
    14 = 2
    15 = 2
    16 = 3
}
foreach ($row in $counts.Keys) {
    $wsParallel.Range("H$row").Value = $counts[$row]
}

# Size the new column similarly to the other label/value columns.
$wsParallel.Columns.Item(8).ColumnWidth = 27

# --- 3. Shift the big comparison picture one column over (col J -> K) ---
foreach ($shape in $wsParallel.Shapes) {
    if ($shape.Name -eq "Picture 13") {
        $shape.Left = $wsParallel.Columns.Item(11).Left
    }
}

# --- 4. Print in portrait orientation ---
$wsParallel.PageSetup.Orientation = 1

# --- 5. Scroll the Perpendicular_185859 sheet down a bit ---
$wsPerpendicular.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1

# --- 6. Make Parallel_185859 the active sheet/selection ---
$wsParallel.Activate()
$wsParallel.Range("L13").Select()
